# Actualizacion automatica 2025-06-05 10:49:05
# Adds the "CUMPLIMIENTO MENSUAL" sheet (a per-ASESOR / per-GRUPO compliance
# summary: PRESUPUESTO, VENTA, POR CUMPLIR, CUMPLIMIENTO) after "VENTA MENSUAL".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new worksheet at the end of the tab strip and rename it.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "CUMPLIMIENTO MENSUAL"

# ---------------------------------------------------------------------------
# 2. Header row (bold, centered, bordered - reuse the same look already used
#    by the header rows on the other two sheets so no superfluous style gets
#    created).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

$headerFormatSource = $wb.Worksheets.Item(1).Range("A1:F1")
$headerFormatSource.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Data rows: ASESOR / GRUPO / PRESUPUESTO / VENTA / POR CUMPLIR / CUMPLIMIENTO
# ---------------------------------------------------------------------------
$data = @(
  ,@("CASTRO ALCIVAR EDA MARIA", "240X120 PORCELANATO", 5820, 1669.25, 4150.75, 0.2868127147766323)
  ,@("CASTRO ALCIVAR EDA MARIA", "240X80 PORCELANATO", 13728, 0, 13728, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "FREGADEROS DE COCINA", 646, 128.74, 517.26, 0.1992879256965945)
  ,@("CASTRO ALCIVAR EDA MARIA", "GRANITO", 238.32, 0, 238.32, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "GRIFERIAS", 106.82, 0, 106.82, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "INODOROS", 2100, 0, 2100, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "LAVABOS", 1000, 0, 1000, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "LED", 300, 0, 300, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "NO RESURTIBLES", 1300.5, 0, 1300.5, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "OTROS", 0, 0, 0, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "PANELES DECORATIVOS", 350, 0, 350, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "PANELES PU", 230, 0, 230, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "PANELES PVC", 966, 0, 966, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "PIEDRA SINTERIZADA", 15690, 2568.3, 13121.7, 0.1636902485659656)
  ,@("CASTRO ALCIVAR EDA MARIA", "PORCELANATO", 45745.689, 346.47, 45399.219, 0.007573828432226697)
  ,@("CASTRO ALCIVAR EDA MARIA", "PUERTAS DE SEGURIDAD", 1142, 0, 1142, 0)
  ,@("CASTRO ALCIVAR EDA MARIA", "SAL SOLUBLE", 1600, 0, 1600, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$lastDataRow = $r - 1
$totalRow = $r

# ---------------------------------------------------------------------------
# 4. TOTAL row.
# ---------------------------------------------------------------------------
$ws.Cells.Item($totalRow, 2).Value = "TOTAL"
$ws.Cells.Item($totalRow, 3).Value = 90963.329
$ws.Cells.Item($totalRow, 4).Value = 4712.76
$ws.Cells.Item($totalRow, 5).Value = 86250.56899999999
$ws.Cells.Item($totalRow, 6).Value = 0.05180944949804993

# ---------------------------------------------------------------------------
# 5. Number formats.
#    C,D,E (PRESUPUESTO / VENTA / POR CUMPLIR) -> same currency format as the
#    other sheets; F (CUMPLIMIENTO) -> percent. Applying to the whole data +
#    total block reuses the existing styles instead of minting new ones.
# ---------------------------------------------------------------------------
$ws.Range("C2:E" + $totalRow).NumberFormat = """$""#,##0.00"
$ws.Range("F2:F" + $totalRow).NumberFormat = "0.00%"

# Right-align the "TOTAL" label (B<totalRow>) - general number format, just
# right aligned, matching the lone new style the diff introduces.
$ws.Cells.Item($totalRow, 2).HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 6. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.1666666666667
$ws.Columns.Item(2).ColumnWidth = 21.1666666666667
$ws.Columns.Item(3).ColumnWidth = 16.1666666666667
$ws.Columns.Item(4).ColumnWidth = 12.1666666666667
$ws.Columns.Item(5).ColumnWidth = 16.1666666666667
$ws.Columns.Item(6).ColumnWidth = 25.1666666666667

# ---------------------------------------------------------------------------
# 7. Keep the original active sheet ("VENTAS POR GRUPO") selected, as the new
#    sheet is only appended - it does not take over the active tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
